$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix a single mislabeled entry: "Cauliflower...123" -> "Cauliflower" (row 118, col B)
$ws.Range("B118").Value = "Cauliflower"

# 2. The questionnaire rows 137-164 (col B) are replaced with the vegetable list that
#    used to live at rows 164-189 (the "milk/fat/salt/supplement" questions are dropped).
$newValues = @(
    "Fish and fish products",
    "Meat, meat products and meat dishes (including bacon, ham and chicken)",
    "Asparagus",
    "Artichoke",
    "Beansprouts...171",
    "Beetroot",
    "Broad beans",
    "Brocoli",
    "Brussel sprouts",
    "Cauliflower...176",
    "Cabbage...177",
    "Chick peas",
    "courgette",
    "Curly kale",
    "Green beans",
    "Leeks...182",
    "Lentils",
    "Lettuce",
    "Mixed Veg Frozen",
    "Mixed Veg Canned",
    "Parsnips",
    "Peas...188",
    "Red Kidney Beans",
    "Runner Beans",
    "Spinach fresh",
    "Spinach frozen",
    "Sweetcorn fresh",
    "Sweetcorn canned"
)

$startRow = 137
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 2).Value = $newValues[$i]
}

# 3. The trailing rows (165-198), which held the now-removed milk/fat/supplement
#    follow-up questions and old vegetable duplicates, are deleted entirely.
$ws.Range("A165:A198").EntireRow.Delete()
